$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 9999
$ws.Range("I21").Value = 9999
$ws.Range("K21").Value = 9999
$ws.Range("M21").Value = -9531

$ws.Range("H23").Value = 9999
$ws.Range("I23").Value = 9999
$ws.Range("K23").Value = 9999
$ws.Range("M23").Value = -9765

$ws.Range("H33").Value = 95331.625
$ws.Range("I33").Value = 151330
$ws.Range("J33").Value = 2001
$ws.Range("K33").Value = 151330
$ws.Range("L33").Value = 2001
$ws.Range("M33").Value = -151101
$ws.Range("N33").Value = -2459

$ws.Range("H86").Value = 5156778.5
$ws.Range("I86").Value = 2540.7058
$ws.Range("K86").Value = 2540.7058
$ws.Range("M86").Value = -1417.7058

$ws.Range("H89").Value = 5156778.5
$ws.Range("I89").Value = 2540.7058
$ws.Range("K89").Value = 2540.7058
$ws.Range("M89").Value = -7087.529

$ws.Range("H111").Value = 20531.223
$ws.Range("I111").Value = 10243.4
$ws.Range("J111").Value = 33391
$ws.Range("K111").Value = 30730.2
$ws.Range("L111").Value = 100173
$ws.Range("M111").Value = -27663.2
$ws.Range("N111").Value = -106307

$ws.Range("H113").Value = 66670564
$ws.Range("I113").Value = 125002380
$ws.Range("K113").Value = 125002380
$ws.Range("M113").Value = -124999126

$ws.Range("H132").Value = 3565.1707
$ws.Range("I132").Value = 3242.3713
$ws.Range("J132").Value = 5448.1665
$ws.Range("K132").Value = 9727.1139
$ws.Range("L132").Value = 16344.4995
$ws.Range("M132").Value = -7197.1139
$ws.Range("N132").Value = -21404.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 24999.8
$ws.Range("J4").Value = 24999.8
$ws.Range("L4").Value = 24999.8
$ws.Range("N4").Value = -25231.8

$ws.Range("H37").Value = 74998.60000000001
$ws.Range("J37").Value = 74998.60000000001
$ws.Range("L37").Value = 74998.60000000001
$ws.Range("N37").Value = -75544.60000000001

$ws.Range("H61").Value = 991.1081
$ws.Range("I61").Value = 991.1081
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 991.1081
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -779.1081
$ws.Range("N61").ClearContents()

$ws.Range("H88").Value = 27779060
$ws.Range("I88").Value = 83334340
$ws.Range("J88").Value = 1421.75
$ws.Range("K88").Value = 83334340
$ws.Range("L88").Value = 1421.75
$ws.Range("M88").Value = -83333934
$ws.Range("N88").Value = -2233.75

$ws.Range("H91").Value = 27779060
$ws.Range("I91").Value = 83334340
$ws.Range("J91").Value = 1421.75
$ws.Range("K91").Value = 83334340
$ws.Range("L91").Value = 1421.75
$ws.Range("M91").Value = -83332936
$ws.Range("N91").Value = -4229.75

$ws.Range("H110").Value = 125007450
$ws.Range("I110").Value = 142864220
$ws.Range("J110").Value = 10000
$ws.Range("K110").Value = 142864220
$ws.Range("L110").Value = 10000
$ws.Range("M110").Value = -142862175
$ws.Range("N110").Value = -14090

$ws.Range("H136").Value = 991.1081
$ws.Range("I136").Value = 991.1081
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2973.3243
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -423.3243000000002
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11630408
$ws.Range("I86").Value = 14708248
$ws.Range("J86").Value = 3011.5557
$ws.Range("K86").Value = 14708248
$ws.Range("L86").Value = 3011.5557
$ws.Range("M86").Value = -14707125
$ws.Range("N86").Value = -5257.5557

$ws.Range("H89").Value = 11630408
$ws.Range("I89").Value = 14708248
$ws.Range("J89").Value = 3011.5557
$ws.Range("K89").Value = 73541240
$ws.Range("L89").Value = 15057.7785
$ws.Range("M89").Value = -73535624
$ws.Range("N89").Value = -26289.7785

$ws.Range("H99").Value = 1836.8334
$ws.Range("I99").Value = 1654.7
$ws.Range("K99").Value = 1654.7
$ws.Range("M99").Value = -156.7

$ws.Range("H105").Value = 1874.75
$ws.Range("I105").Value = 1733
$ws.Range("K105").Value = 1733
$ws.Range("M105").Value = 14

$ws.Range("H134").Value = 1292.5555
$ws.Range("I134").Value = 1292.5555
$ws.Range("K134").Value = 3877.6665
$ws.Range("M134").Value = -1342.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4928.604
$ws.Range("I31").Value = 9127.866
$ws.Range("J31").Value = 3271
$ws.Range("K31").Value = 9127.866
$ws.Range("L31").Value = 3271
$ws.Range("M31").Value = -8832.866
$ws.Range("N31").Value = -3861

$ws.Range("H34").Value = 4928.604
$ws.Range("I34").Value = 9127.866
$ws.Range("J34").Value = 3271
$ws.Range("K34").Value = 9127.866
$ws.Range("L34").Value = 3271
$ws.Range("M34").Value = -8925.866
$ws.Range("N34").Value = -3675

$ws.Range("H58").Value = 1747.7428
$ws.Range("J58").Value = 2273.3
$ws.Range("L58").Value = 2273.3
$ws.Range("N58").Value = -2679.3

$ws.Range("H136").Value = 1747.7428
$ws.Range("J136").Value = 2273.3
$ws.Range("L136").Value = 6819.900000000001
$ws.Range("N136").Value = -11919.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11153643
$ws.Range("I4").Value = 18033068
$ws.Range("J4").Value = 4231.4136
$ws.Range("K4").Value = 54099204
$ws.Range("L4").Value = 12694.2408
$ws.Range("M4").Value = -54099092
$ws.Range("N4").Value = -12918.2408

$ws.Range("H5").Value = 1130.0952
$ws.Range("I5").Value = 656.6667
$ws.Range("J5").Value = 1485.1666
$ws.Range("K5").Value = 1970.0001
$ws.Range("L5").Value = 4455.4998
$ws.Range("M5").Value = -1858.0001
$ws.Range("N5").Value = -4679.4998

$ws.Range("H86").Value = 874.5
$ws.Range("J86").Value = 999.5
$ws.Range("L86").Value = 2998.5
$ws.Range("N86").Value = -5370.5

$ws.Range("H89").Value = 874.5
$ws.Range("J89").Value = 999.5
$ws.Range("L89").Value = 8995.5
$ws.Range("N89").Value = -20851.5

$ws.Range("H113").Value = 1364.1364
$ws.Range("J113").Value = 1772.3846
$ws.Range("L113").Value = 5317.1538
$ws.Range("N113").Value = -9657.1538

$ws.Range("H131").Value = 5156.696
$ws.Range("J131").Value = 7445.8
$ws.Range("L131").Value = 22337.4
$ws.Range("N131").Value = -32417.4

$ws.Range("H135").Value = 1130.0952
$ws.Range("I135").Value = 656.6667
$ws.Range("J135").Value = 1485.1666
$ws.Range("K135").Value = 5910.0003
$ws.Range("L135").Value = 13366.4994
$ws.Range("M135").Value = -3375.0003
$ws.Range("N135").Value = -18436.4994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11487.5
$ws.Range("I70").Value = 12666.667
$ws.Range("K70").Value = 12666.667
$ws.Range("M70").Value = -12396.667

$ws.Range("H73").Value = 11487.5
$ws.Range("I73").Value = 12666.667
$ws.Range("K73").Value = 12666.667
$ws.Range("M73").Value = -11730.667

$ws.Range("H80").Value = 3754.923
$ws.Range("I80").Value = 3707.4707
$ws.Range("J80").Value = 3844.5557
$ws.Range("K80").Value = 3707.4707
$ws.Range("L80").Value = 3844.5557
$ws.Range("M80").Value = -2709.4707
$ws.Range("N80").Value = -5840.5557

$ws.Range("H83").Value = 3754.923
$ws.Range("I83").Value = 3707.4707
$ws.Range("J83").Value = 3844.5557
$ws.Range("K83").Value = 18537.3535
$ws.Range("L83").Value = 19222.7785
$ws.Range("M83").Value = -13545.3535
$ws.Range("N83").Value = -29206.7785

$ws.Range("H123").Value = 32475
$ws.Range("J123").Value = 32475
$ws.Range("L123").Value = 32475
$ws.Range("N123").Value = -37375

$ws.Range("H132").Value = 21380.44
$ws.Range("I132").Value = 29096.625
$ws.Range("K132").Value = 87289.875
$ws.Range("M132").Value = -84759.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2116.0322
$ws.Range("I46").Value = 1311.7646
$ws.Range("J46").Value = 3092.6428
$ws.Range("K46").Value = 1311.7646
$ws.Range("L46").Value = 3092.6428
$ws.Range("M46").Value = -1123.7646
$ws.Range("N46").Value = -3468.6428

$ws.Range("H100").Value = 2052.5
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459

$ws.Range("H132").Value = 6094.577
$ws.Range("I132").Value = 3562.5625
$ws.Range("K132").Value = 10687.6875
$ws.Range("M132").Value = -8157.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1054761.4
$ws.Range("I62").Value = 2651502.8
$ws.Range("J62").Value = 28284.715
$ws.Range("K62").Value = 2651502.8
$ws.Range("L62").Value = 28284.715
$ws.Range("M62").Value = -2650878.8
$ws.Range("N62").Value = -29532.715

$ws.Range("H65").Value = 1054761.4
$ws.Range("I65").Value = 2651502.8
$ws.Range("J65").Value = 28284.715
$ws.Range("K65").Value = 13257514
$ws.Range("L65").Value = 141423.575
$ws.Range("M65").Value = -13254394
$ws.Range("N65").Value = -147663.575

$ws.Range("H100").Value = 10722
$ws.Range("I100").Value = 10722
$ws.Range("K100").Value = 21444
$ws.Range("M100").Value = -20903

$ws.Range("H132").Value = 2584.7544
$ws.Range("I132").Value = 2361.348
$ws.Range("J132").Value = 3519
$ws.Range("K132").Value = 7084.044
$ws.Range("L132").Value = 10557
$ws.Range("M132").Value = -4554.044
$ws.Range("N132").Value = -15617
